# Apply the "changes to plots of sensitivity analysis" edit:
#  - Rename Sheet1 -> "parameter values", Sheet2 -> "scenario tracker"
#  - On "scenario tracker": every status in column B (rows 3-16) becomes "output .csvs"
#    (this also frees up the now-unused "running"/"added to db"/"run" shared
#    strings, which shifts later shared-string indices, e.g. the reference
#    used by "parameter values"!H45)
#  - "scenario tracker" becomes the active/selected sheet, with B15 selected
#  - "parameter values" keeps its previous selection (H45), but is no longer the active tab

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("Sheet1")
$wsTracker = $wb.Worksheets.Item("Sheet2")

$wsParams.Name = "parameter values"
$wsTracker.Name = "scenario tracker"

# Update the status column on the scenario tracker sheet: all scenarios are
# now finished and their output .csvs have been produced.
$wsTracker.Range("B3:B16").Value = "output .csvs"

# Make "scenario tracker" the active sheet/tab and select cell B15, matching
# the new saved view state.
$wsTracker.Activate() | Out-Null
$wsTracker.Range("B15").Select() | Out-Null
